$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks numeric need to be force-written as text
# (matching the source workbook, where every Price/Volume cell is an
# inline string) - otherwise Excel COM auto-converts "229.42" etc. to a
# real number. We flip NumberFormat to "@" for the write, then restore
# the cell to the default "Normal" style so no stray formatting sticks.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "37.946.29"
$ws.Range("E2").Value = "  +1.99%  "

$ws.Range("D3").Value = "2.050.06"
$ws.Range("E3").Value = "  +1.08%  "

$ws.Range("E4").Value = "  -0.05%  "

Set-TextValue $ws.Range("D5") "229.42"
$ws.Range("E5").Value = "  +1.09%  "

$ws.Range("E6").Value = "  +1.12%  "

Set-TextValue $ws.Range("D7") "58.47"
$ws.Range("E7").Value = "  +5.92%  "

$ws.Range("E8").Value = "  -0.01%  "

Set-TextValue $ws.Range("D9") "0.385"
$ws.Range("E9").Value = "  +1.42%  "

$ws.Range("E10").Value = "  +2.91%  "

$ws.Range("E11").Value = "  +1.37%  "

$ws.Range("D12").Value = "2.353.27"
$ws.Range("E12").Value = "  +1.36%  "

$ws.Range("E13").Value = "  +2.45%  "

Set-TextValue $ws.Range("D14") "20.83"
$ws.Range("E14").Value = "  +2.57%  "

Set-TextValue $ws.Range("D15") "0.752"
$ws.Range("E15").Value = "  +1.05%  "

Set-TextValue $ws.Range("D16") "5.28"
$ws.Range("E16").Value = "  +1.69%  "

$ws.Range("D17").Value = "2.054.83"
$ws.Range("E17").Value = "  -0.06%  "

$ws.Range("D18").Value = "37.880.65"

Set-TextValue $ws.Range("D19") "6.25"
$ws.Range("E19").Value = "  -3.12%  "

Set-TextValue $ws.Range("D20") "69.70"
$ws.Range("E20").Value = "  +1.07%  "

$ws.Range("E21").Value = "  +2.20%  "

Set-TextValue $ws.Range("D22") "224.72"
$ws.Range("E22").Value = "  +0.33%  "

$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("E25").Value = "  +2.27%  "

$ws.Range("E26").Value = "  +0.26%  "

Set-TextValue $ws.Range("D27") "166.41"
$ws.Range("E27").Value = "  +0.12%  "

$ws.Range("E28").Value = "  +4.39%  "

Set-TextValue $ws.Range("D29") "19.00"
$ws.Range("E29").Value = "  +1.26%  "

$ws.Range("E30").Value = "  +0.36%  "

$ws.Range("E31").Value = "  +1.16%  "

$ws.Range("E32").Value = "  -0.32%  "

Set-TextValue $ws.Range("D33") "4.58"
$ws.Range("E33").Value = "  +2.35%  "

$ws.Range("E34").Value = "  -0.61%  "

$ws.Range("E35").Value = "  +9.28%  "

$ws.Range("E36").Value = "  -1.19%  "

$ws.Range("E37").Value = "  +9.36%  "

$ws.Range("E38").Value = "  +5.36%  "

$ws.Range("E39").Value = "  +0.03%  "

$ws.Range("D40").Value = "1.483.86"
$ws.Range("E40").Value = "  +0.74%  "

$ws.Range("E41").Value = "  +0.42%  "

Set-TextValue $ws.Range("D42") "97.18"
$ws.Range("E42").Value = "  +1.24%  "

$ws.Range("E43").Value = "  +3.48%  "

Set-TextValue $ws.Range("D44") "16.49"
$ws.Range("E44").Value = "  +0.61%  "

$ws.Range("E45").Value = "  +0.93%  "

$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D46") "1.13"
$ws.Range("E46").Value = "  -1.36%  "

$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D47") "4.14"
$ws.Range("E47").Value = "  +13.05%  "

$ws.Range("E48").Value = "  +0.16%  "

$ws.Range("E49").Value = "  +1.24%  "

$ws.Range("E50").Value = "  -2.83%  "

$ws.Range("D51").Value = "2.242.30"
$ws.Range("E51").Value = "  +1.59%  "
